# Code changes for THS template E!
$wb = $excel.ActiveWorkbook

$wsTestData = $wb.Worksheets.Item("TestData")
$wsWindows  = $wb.Worksheets.Item("Windows")

# --- TestData sheet: add new TC3 row (row 4), rename existing TC2 row's TcNo to TC3 ---

# Insert a new row 4, copying formatting down from row 3 (this is what gives the
# new row's cells the same "touched" style index (s="4") as columns B/C in row 3).
[void]$wsTestData.Rows.Item(4).Insert(-4121, 0)

# Row 3 had style index 4 on every column; the freshly inserted row 4 copied that
# onto every column too. Per the target, only columns B and C (Distributor, DealType)
# keep that style on row 4 - clear the formatting on A4/D4/E4 back to default.
[void]$wsTestData.Range("A4").ClearFormats()
[void]$wsTestData.Range("D4").ClearFormats()
[void]$wsTestData.Range("E4").ClearFormats()

# The existing TC2 row (row 2) becomes TC3
$wsTestData.Range("A2").Value = "TC3"

# New row 4 holds the old TC2 data plus a new distributor/negotiator pairing
$wsTestData.Range("A4").Value = "TC2"
$wsTestData.Range("B4").Value = "10X10 Entertainment"
$wsTestData.Range("C4").Value = "Cash"
$wsTestData.Range("E4").Value = "Esquire Network"
$wsTestData.Range("D4").Value = "Doug Baughman"

# --- Windows sheet: swap TC1/TC2 labels on rows 2-4, update last window's dates ---

$wsWindows.Range("A2").Value = "TC2"
$wsWindows.Range("A3").Value = "TC2"
$wsWindows.Range("A4").Value = "TC1"

$wsWindows.Range("C4").Value = [datetime]"2018-12-01"
$wsWindows.Range("D4").Value = [datetime]"2018-12-31"

# --- Active sheet / selection bookkeeping ---
# TestData is no longer the selected tab; its last selection moves to A4.
[void]$wsTestData.Activate()
[void]$wsTestData.Range("A4").Select()

# Windows becomes the active/selected tab, with D5 as its last selection.
[void]$wsWindows.Activate()
[void]$wsWindows.Range("D5").Select()
